# Regenerate save_data: column G (header "K", formerly "Strike#"-derived) is
# recomputed from the regenerated std/mean "s_vals" calculation. This writes
# the freshly-calculated K values back onto each data row (rows 2-84) of the
# active sheet, leaving every other column untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new value for column G ("K").
$kValues = @{
    2  = 0
    3  = 0
    4  = 0
    5  = 1
    6  = 0
    7  = 1
    8  = 1
    9  = 0
    10 = 1
    11 = 1
    12 = 1
    13 = 1
    14 = 0
    15 = 0
    16 = 2
    17 = 0
    18 = 2
    19 = 0
    20 = 1
    21 = 4
    22 = 0
    23 = 0
    24 = 0
    25 = 0
    26 = 0
    27 = 1
    28 = 2
    29 = 0
    30 = 0
    31 = 2
    32 = 0
    33 = 1
    34 = 0
    35 = 2
    36 = 3
    37 = 2
    38 = 0
    39 = 1
    40 = 3
    41 = 4
    42 = 0
    43 = 1
    44 = 0
    45 = 1
    46 = 1
    47 = 0
    48 = 2
    49 = 3
    50 = 2
    51 = 1
    52 = 0
    53 = 1
    54 = 1
    55 = 0
    56 = 1
    57 = 1
    58 = 0
    59 = 1
    60 = 0
    61 = 1
    62 = 1
    63 = 0
    64 = 0
    65 = 2
    66 = 2
    67 = 3
    68 = 0
    69 = 4
    70 = 0
    71 = 1
    72 = 1
    73 = 2
    74 = 2
    75 = 0
    76 = 0
    77 = 1
    78 = 2
    79 = 0
    80 = 1
    81 = 1
    82 = 2
    83 = 1
    84 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
